$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2623293280092867
$ws.Range("C2").Value = 0.04519528119864447
$ws.Range("D2").Value = 0.07837785008241838
$ws.Range("E2").Value = 0.1560924885874613
$ws.Range("G2").Value = 0.002453381699923774
$ws.Range("I2").Value = 0.7232617069412512
$ws.Range("K2").Value = 0.2856344520822347
$ws.Range("M2").Value = 0.2290699548547721
$ws.Range("O2").Value = 3.387746781384607
$ws.Range("B3").Value = 0.2323202765363135
$ws.Range("C3").Value = 0.03941117441554809
$ws.Range("D3").Value = 0.07111882269654757
$ws.Range("E3").Value = 0.1451477278358695
$ws.Range("G3").Value = 0.00245619341695368
$ws.Range("I3").Value = 0.7262354139395875
$ws.Range("K3").Value = 0.2513121107416794
$ws.Range("M3").Value = 0.2070729657989219
$ws.Range("O3").Value = 3.387659030955064
$ws.Range("B4").Value = 0.2139124479455745
$ws.Range("C4").Value = 0.03584950314488822
$ws.Range("D4").Value = 0.06669625504463284
$ws.Range("E4").Value = 0.1385324210059835
$ws.Range("G4").Value = 0.002458011640592335
$ws.Range("I4").Value = 0.728436962065679
$ws.Range("K4").Value = 0.2302386408494357
$ws.Range("M4").Value = 0.1936562712843113
$ws.Range("O4").Value = 3.389492456128721
$ws.Range("B5").Value = 0.2064159351005515
$ws.Range("C5").Value = 0.03439554240182474
$ws.Range("D5").Value = 0.06490269447211006
$ws.Range("E5").Value = 0.135862742392554
$ws.Range("G5").Value = 0.002458775740481219
$ws.Range("I5").Value = 0.7294284985647508
$ws.Range("K5").Value = 0.2216514970477732
$ws.Range("M5").Value = 0.1882112737812065
$ws.Range("O5").Value = 3.390713729406514
$ws.Range("B6").Value = 0.2051714463941323
$ws.Range("C6").Value = 0.0341539597431364
$ws.Range("D6").Value = 0.06460539872111326
$ws.Range("E6").Value = 0.1354210142739305
$ws.Range("G6").Value = 0.002458904019492745
$ws.Range("I6").Value = 0.7295988408983867
$ws.Range("K6").Value = 0.2202256459826515
$ws.Range("M6").Value = 0.1873084869241595
$ws.Range("O6").Value = 3.3909451448217
$ws.Range("B7").Value = 0.2138113273578313
$ws.Range("C7").Value = 0.03582990483009496
$ws.Range("D7").Value = 0.06667203135367572
$ws.Range("E7").Value = 0.1384963113808908
$ws.Range("G7").Value = 0.002458021851520513
$ws.Range("I7").Value = 0.7284499522096795
$ws.Range("K7").Value = 0.2301228293123216
$ws.Range("M7").Value = 0.1935827475220435
$ws.Range("O7").Value = 3.389507007477562
$ws.Range("B8").Value = 0.2519786950874163
$ws.Range("C8").Value = 0.04320304707904654
$ws.Range("D8").Value = 0.0758677708615636
$ws.Range("E8").Value = 0.1522968283570165
$ws.Range("G8").Value = 0.002454332165605622
$ws.Range("I8").Value = 0.7242090196113047
$ws.Range("K8").Value = 0.2738001670693677
$ws.Range("M8").Value = 0.2214667389026062
$ws.Range("O8").Value = 3.387324540732749
$ws.Range("B9").Value = 0.326955387947919
$ws.Range("C9").Value = 0.05758080168980939
$ws.Range("D9").Value = 0.09417555816767731
$ws.Range("E9").Value = 0.18020383501986
$ws.Range("G9").Value = 0.002447822004554855
$ws.Range("I9").Value = 0.7188775271108554
$ws.Range("K9").Value = 0.35944627564146
$ws.Range("M9").Value = 0.2768652335752364
$ws.Range("O9").Value = 3.398043876733141
$ws.Range("B10").Value = 0.3821111932499832
$ws.Range("C10").Value = 0.06809588302058955
$ws.Range("D10").Value = 0.1077969756150878
$ws.Range("E10").Value = 0.2012412366257905
$ws.Range("G10").Value = 0.002443476619471668
$ws.Range("I10").Value = 0.7167866100810869
$ws.Range("K10").Value = 0.4223608153706948
$ws.Range("M10").Value = 0.318018978519568
$ws.Range("O10").Value = 3.415103115468639
$ws.Range("B11").Value = 0.4072167209222073
$ws.Range("C11").Value = 0.07286929255130303
$ws.Range("D11").Value = 0.1140315090632242
$ws.Range("E11").Value = 0.210931992321477
$ws.Range("G11").Value = 0.002441593853526144
$ws.Range("I11").Value = 0.7162333618086834
$ws.Range("K11").Value = 0.4509792023866055
$ws.Range("M11").Value = 0.3368425603689076
$ws.Range("O11").Value = 3.424866981164371
$ws.Range("B12").Value = 0.4167254141395347
$ws.Range("C12").Value = 0.07467542878538325
$ws.Range("D12").Value = 0.1163978675011776
$ws.Range("E12").Value = 0.2146192940195846
$ws.Range("G12").Value = 0.002440894339056986
$ws.Range("I12").Value = 0.7160811877767372
$ws.Range("K12").Value = 0.4618157652849106
$ws.Range("M12").Value = 0.3439854891673733
$ws.Range("O12").Value = 3.428853025625301
$ws.Range("B13").Value = 0.4146774737658347
$ws.Range("C13").Value = 0.07428650983523255
$ws.Range("D13").Value = 0.1158879863364746
$ws.Range("E13").Value = 0.2138243796785133
$ws.Range("G13").Value = 0.002441044394749793
$ws.Range("I13").Value = 0.7161114097195878
$ws.Range("K13").Value = 0.4594819497226013
$ws.Range("M13").Value = 0.3424464696919713
$ws.Range("O13").Value = 3.427981711328698
$ws.Range("B14").Value = 0.4079989733461105
$ws.Range("C14").Value = 0.07301791364153587
$ws.Range("D14").Value = 0.1142260810531752
$ws.Range("E14").Value = 0.2112349941649114
$ws.Range("G14").Value = 0.002441536034680014
$ws.Range("I14").Value = 0.7162196927768534
$ws.Range("K14").Value = 0.4518707469259482
$ws.Range("M14").Value = 0.3374299152779656
$ws.Range("O14").Value = 3.42518912656115
$ws.Range("B15").Value = 0.4039084193539964
$ws.Range("C15").Value = 0.07224067240595389
$ws.Range("D15").Value = 0.11320882919388
$ws.Range("E15").Value = 0.2096512235471764
$ws.Range("G15").Value = 0.002441838927802929
$ws.Range("I15").Value = 0.7162934884400229
$ws.Range("K15").Value = 0.4472085771552941
$ws.Range("M15").Value = 0.3343590673993049
$ws.Range("O15").Value = 3.423516199539222
$ws.Range("B16").Value = 0.3804707561066891
$ws.Range("C16").Value = 0.06778372755168505
$ws.Range("D16").Value = 0.1073903004924404
$ws.Range("E16").Value = 0.2006103747751027
$ws.Range("G16").Value = 0.002443601547980423
$ws.Range("I16").Value = 0.7168307826141955
$ws.Range("K16").Value = 0.4204904714756879
$ws.Range("M16").Value = 0.3167908851976549
$ws.Range("O16").Value = 3.414505389209921
$ws.Range("B17").Value = 0.366096058800963
$ws.Range("C17").Value = 0.06504697319648756
$ws.Range("D17").Value = 0.1038305753312301
$ws.Range("E17").Value = 0.19509521587009
$ws.Range("G17").Value = 0.002444706880456726
$ws.Range("I17").Value = 0.7172623893931203
$ws.Range("K17").Value = 0.4040990823569643
$ws.Range("M17").Value = 0.3060397090479654
$ws.Range("O17").Value = 3.40949111057887
$ws.Range("B18").Value = 0.3578295348597749
$ws.Range("C18").Value = 0.06347192869685614
$ws.Range("D18").Value = 0.1017867001017976
$ws.Range("E18").Value = 0.1919343909507916
$ws.Range("G18").Value = 0.002445351487152629
$ws.Range("I18").Value = 0.7175480858371515
$ws.Range("K18").Value = 0.3946710444976702
$ws.Range("M18").Value = 0.2998655619844186
$ws.Range("O18").Value = 3.406795577652076
$ws.Range("B19").Value = 0.3550308872169126
$ws.Range("C19").Value = 0.06293848555918657
$ws.Range("D19").Value = 0.1010952947158898
$ws.Range("E19").Value = 0.1908661311331841
$ws.Range("G19").Value = 0.002445571261673342
$ws.Range("I19").Value = 0.7176512461081614
$ws.Range("K19").Value = 0.3914788586949385
$ws.Range("M19").Value = 0.2977767553281652
$ws.Range("O19").Value = 3.405915282243399
$ws.Range("B20").Value = 0.3676261253701512
$ws.Range("C20").Value = 0.06533840235867672
$ws.Range("D20").Value = 0.1042091431671111
$ws.Range("E20").Value = 0.1956811375246659
$ws.Range("G20").Value = 0.002444588300458207
$ws.Range("I20").Value = 0.717212567635471
$ws.Range("K20").Value = 0.4058439911769938
$ws.Range("M20").Value = 0.3071831911358558
$ws.Range("O20").Value = 3.410005370963944
$ws.Range("B21").Value = 0.4099605655334244
$ws.Range("C21").Value = 0.073390570606648
$ws.Range("D21").Value = 0.1147140743087363
$ws.Range("E21").Value = 0.2119950790945637
$ws.Range("G21").Value = 0.002441391263739945
$ws.Range("I21").Value = 0.7161863306570098
$ws.Range("K21").Value = 0.4541063607009903
$ws.Range("M21").Value = 0.3389029951537665
$ws.Range("O21").Value = 3.426001537111858
$ws.Range("B22").Value = 0.4376387466494123
$ws.Range("C22").Value = 0.07864466569134265
$ws.Range("D22").Value = 0.1216115715791091
$ws.Range("E22").Value = 0.2227600240859573
$ws.Range("G22").Value = 0.002439380171776473
$ws.Range("I22").Value = 0.7158498337229915
$ws.Range("K22").Value = 0.4856449574154169
$ws.Range("M22").Value = 0.3597203077544791
$ws.Range("O22").Value = 3.43813887404005
$ws.Range("B23").Value = 0.4228655690270102
$ws.Range("C23").Value = 0.07584123675019327
$ws.Range("D23").Value = 0.1179273272096708
$ws.Range("E23").Value = 0.2170050738974183
$ws.Range("G23").Value = 0.00244044638115606
$ws.Range("I23").Value = 0.7159988124762862
$ws.Range("K23").Value = 0.4688126629763474
$ws.Range("M23").Value = 0.3486017595210811
$ws.Range("O23").Value = 3.431506771743813
$ws.Range("B24").Value = 0.3669343894720498
$ws.Range("C24").Value = 0.06520665236665479
$ws.Range("D24").Value = 0.1040379843699952
$ws.Range("E24").Value = 0.1954162114886557
$ws.Range("G24").Value = 0.00244464188194381
$ws.Range("I24").Value = 0.7172349750594051
$ws.Range("K24").Value = 0.4050551316118174
$ws.Range("M24").Value = 0.306666201567559
$ws.Range("O24").Value = 3.409772290619969
$ws.Range("B25").Value = 0.3066591835933536
$ws.Range("C25").Value = 0.05369979305244499
$ws.Range("D25").Value = 0.08919306190315979
$ws.Range("E25").Value = 0.172561860710978
$ws.Range("G25").Value = 0.002449505997477757
$ws.Range("I25").Value = 0.7199996136101703
$ws.Range("K25").Value = 0.3362779632364834
$ws.Range("M25").Value = 0.2618000806283334
$ws.Range("O25").Value = 3.393534133380058
